# Update column F ("想去人数") values across the four worksheets to match
# the regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$sheet1 = @{
    2  = 76
    3  = 279
    4  = 632
    5  = 2794
    6  = 14
    7  = 221
    8  = 272
    9  = 6451
    12 = 77
    13 = 5050
    14 = 111
    15 = 1
    17 = 2647
    18 = 1358
    19 = 1524
    20 = 1230
    21 = 325
    22 = 124
    24 = 1086
    25 = 250
    26 = 1392
    27 = 1058
    28 = 2117
    29 = 330
    30 = 591
    31 = 37
    32 = 38
    33 = 109
    34 = 256
    35 = 1521
    37 = 5
    38 = 616
    39 = 1080
    40 = 120
    42 = 314
    43 = 2310
    44 = 2573
    46 = 140
    47 = 103
    48 = 110
}
foreach ($row in $sheet1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1[$row]
}

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$sheet2 = @{
    8  = 1
    17 = 160
    21 = 153
    24 = 5
    27 = 416
    41 = 29
}
foreach ($row in $sheet2.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2[$row]
}

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$sheet3 = @{
    5  = 21
    6  = 1707
    7  = 574
    8  = 1529
    9  = 1821
    10 = 2568
    11 = 875
    12 = 750
    14 = 133
}
foreach ($row in $sheet3.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3[$row]
}

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$sheet4 = @{
    2  = 76
    4  = 1707
    5  = 632
    6  = 574
    7  = 2794
    8  = 221
    9  = 1529
    10 = 272
    11 = 2568
    12 = 6452
    13 = 750
    16 = 1358
    17 = 1524
    18 = 325
    19 = 124
    21 = 1086
    22 = 250
    24 = 1392
    25 = 1058
    26 = 2117
    27 = 330
    28 = 591
    29 = 37
    30 = 38
    31 = 256
    33 = 1521
    35 = 5
    36 = 1080
    39 = 416
    40 = 314
    43 = 2310
    44 = 2573
    45 = 103
    46 = 110
}
foreach ($row in $sheet4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4[$row]
}
